$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format D cells whose new values look numeric, so they stay text
$textCells = @("D5","D6","D7","D8","D10","D11","D12","D14","D15","D16","D17","D18","D19","D21","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.153.10'
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").Value = '1.904.11'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("D5").Value = '306.21'
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").Value = '0.5223'
$ws.Range("E7").Value = '  +1.54%  '

$ws.Range("D8").Value = '0.3758'
$ws.Range("E8").Value = '  +0.38%  '

$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").Value = '21.12'

$ws.Range("D11").Value = '0.9024'
$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("D12").Value = '0.08520'
$ws.Range("E12").Value = '  +11.52%  '

$ws.Range("D13").Value = '1.910.18'
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").Value = '95.06'
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").Value = '5.290'
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("D17").Value = '0.000008630'
$ws.Range("E17").Value = '  +1.47%  '

$ws.Range("D18").Value = '14.55'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").Value = '27.186.35'
$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").Value = '5.068'
$ws.Range("E21").Value = '  -0.19%  '

$ws.Range("D22").Value = '2.151.24'
$ws.Range("E22").Value = '  +0.76%  '

$ws.Range("D23").Value = '10.61'
$ws.Range("E23").Value = '  +0.44%  '

$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").Value = '2.288'
$ws.Range("E25").Value = '  +3.36%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '146.95'
$ws.Range("E26").Value = '  +0.60%  '

$ws.Range("D27").Value = '1.753'
$ws.Range("E27").Value = '  -2.26%  '

$ws.Range("D28").Value = '18.21'

$ws.Range("D29").Value = '114.98'
$ws.Range("E29").Value = '  +0.33%  '

$ws.Range("D30").Value = '4.811'
$ws.Range("E30").Value = '  -1.00%  '

$ws.Range("D31").Value = '4.901'
$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("D32").Value = '0.09258'
$ws.Range("E32").Value = '  +0.68%  '

$ws.Range("D33").Value = '0.8058'
$ws.Range("E33").Value = '  +4.83%  '

$ws.Range("D34").Value = '0.05055'
$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("D35").Value = '1.235'
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").Value = '3.454'
$ws.Range("E36").Value = '  +4.79%  '

$ws.Range("D37").Value = '2.956'
$ws.Range("E37").Value = '  -0.62%  '

$ws.Range("D38").Value = '2.615'
$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("D39").Value = '0.5713'
$ws.Range("E39").Value = '  +1.95%  '

$ws.Range("D40").Value = '0.01996'
$ws.Range("E40").Value = '  -0.22%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").Value = '9.070'
$ws.Range("E42").Value = '  +1.29%  '

$ws.Range("D43").Value = '6.635'
$ws.Range("E43").Value = '  -0.27%  '

$ws.Range("D44").Value = '116.10'
$ws.Range("E44").Value = '  -1.59%  '

$ws.Range("D45").Value = '0.1518'
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("D46").Value = '0.4868'
$ws.Range("E46").Value = '  +1.29%  '

$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("E48").Value = '  -1.03%  '

$ws.Range("D49").Value = '1.615'
$ws.Range("E49").Value = '  +1.40%  '

$ws.Range("D50").Value = '37.51'
$ws.Range("E50").Value = '  +0.00%  '

$ws.Range("D51").Value = '64.00'
$ws.Range("E51").Value = '  +0.10%  '
